$wb = $excel.ActiveWorkbook
$dbd = $wb.Worksheets.Item("DBD")
$dbs = $wb.Worksheets.Item("DBS")

# ------------------------------------------------------------------
# 1. Make room for the old "CreateEmpNo / CreateDate / LastUpdateEmpNo /
#    LastUpdate" block by inserting 4 fresh rows right after it (it
#    currently occupies rows 23-26, followed by 3 spare blank rows
#    27-29). We then copy that block (values + formats) down into the
#    freshly inserted rows so it keeps looking exactly as it did.
# ------------------------------------------------------------------
$dbd.Rows("28:31").Insert()

$dbd.Range("A23:G26").Copy()
$dbd.Range("A28").PasteSpecial(-4104)   # xlPasteAll

# Renumber the SEQ column (A) for the relocated block.
$dbd.Range("A28").Value = 20
$dbd.Range("A29").Value = 21
$dbd.Range("A30").Value = 22
$dbd.Range("A31").Value = 23

# The 2 now-unused spare blank rows (pushed down to 32-33) are removed.
$dbd.Rows("32:33").Delete()

# ------------------------------------------------------------------
# 2. Re-purpose old rows 23-27 (4 filled + 1 spare blank row) in place
#    for the 5 new balance-check fields. Column G (and the row itself)
#    is left untouched; only B:F are retyped.
# ------------------------------------------------------------------

# Pull matching B/C/D/E formatting (left/right aligned, bordered,
# 標楷體 font) from row 9, which already carries the exact look the
# new rows need.
$dbd.Range("B9:E9").Copy()
$dbd.Range("B23:E23").PasteSpecial(-4122)   # xlPasteFormats
$dbd.Range("B24:E24").PasteSpecial(-4122)
$dbd.Range("B25:E25").PasteSpecial(-4122)
$dbd.Range("B26:E26").PasteSpecial(-4122)
$dbd.Range("B27:E27").PasteSpecial(-4122)
$dbd.Range("E9:E9").Copy()
$dbd.Range("F23").PasteSpecial(-4122)
$dbd.Range("F24").PasteSpecial(-4122)
$dbd.Range("F25").PasteSpecial(-4122)
$dbd.Range("F26").PasteSpecial(-4122)
$dbd.Range("F27").PasteSpecial(-4122)

$dbd.Range("A23").Value = 15
$dbd.Range("B23").Value = "YdBal"
$dbd.Range("C23").Value = "前日餘額"
$dbd.Range("D23").Value = "DECIMAL"
$dbd.Range("E23").Value = 18
$dbd.Range("F23").Value = 2

$dbd.Range("A24").Value = 16
$dbd.Range("B24").Value = "DbAmt"
$dbd.Range("C24").Value = "借方金額"
$dbd.Range("D24").Value = "DECIMAL"
$dbd.Range("E24").Value = 18
$dbd.Range("F24").Value = 2

$dbd.Range("A25").Value = 17
$dbd.Range("B25").Value = "CrAmt"
$dbd.Range("C25").Value = "貸方金額"
$dbd.Range("D25").Value = "DECIMAL"
$dbd.Range("E25").Value = 18
$dbd.Range("F25").Value = 2

$dbd.Range("A26").Value = 18
$dbd.Range("B26").Value = "CoreDbAmt"
$dbd.Range("C26").Value = "核心借方金額"
$dbd.Range("D26").Value = "DECIMAL"
$dbd.Range("E26").Value = 18
$dbd.Range("F26").Value = 2

$dbd.Range("A27").Value = 19
$dbd.Range("B27").Value = "CoreCrAmt"
$dbd.Range("C27").Value = "核心貸方金額"
$dbd.Range("D27").Value = "DECIMAL"
$dbd.Range("E27").Value = 18
$dbd.Range("F27").Value = 2

# B23 keeps the formatting it always had (it was never reformatted),
# so restore its original look after the B9:E9 format paste touched it.
$dbd.Range("B28").Copy()
$dbd.Range("B23").PasteSpecial(-4122)
$dbd.Range("B23").Value = "YdBal"

# ------------------------------------------------------------------
# 3. View state: DBD becomes the active sheet/tab, zoomed to 115%,
#    scrolled near the bottom, with the last edited cell selected.
#    DBS stops being the active tab and scrolls one column to the right.
# ------------------------------------------------------------------
$dbd.Activate()
$excel.ActiveWindow.Zoom = 115
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
$dbd.Range("G31").Select()

$dbs.Activate()
$excel.ActiveWindow.ScrollColumn = 2
$dbd.Activate()
$dbd.Range("G31").Select()
